$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $prompt = $ws.Cells.Item($r, 2).Value()
    if ($prompt -match "(?s)## پرسش\r?\n(.*?)\r?\n\s*## گزینه ها") {
        $question = $matches[1].Trim()
        if ($question.StartsWith("'")) {
            $question = "'" + $question
        }
        $ws.Cells.Item($r, 1).Value = $question
    }
}
